$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.7222044307299725, 0.4737940993224269, 1.438000664540498, 3.058035597739332)"
$ws.Range("C2").Value = "JSU(-1.2327934495956392, 1.4316023314512387, 1.9726011866188928, 6.916957536762931)"
$ws.Range("D2").Value = "NIG(0.9232211583265119, 0.6749433853365115, 1.6768657734140184, 2.712771356539606)"
$ws.Range("E2").Value = "NIG(1.7110035190297417, 1.2687056714422946, 3.338197060092456, 5.938165140844333)"
